# Correcting market share tab for updated scenario 3s
#
# The "years active" flags on the MarketShare sheet were shifted one
# row too low (row 3, columns L:Z) and need to live on row 2 instead
# (columns L:Z), matching the "Albendazole" market-share row rather
# than the "Mebendazole" row below it. Also re-point the active tab /
# selection at the MarketShare sheet, on the corrected range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MarketShare")

# Move the stray 1's from row 3 (L3:Z3) up to row 2 (L2:Z2).
$ws.Range("L2:Z2").Value = 1
$ws.Range("L3:Z3").ClearContents()

# Make MarketShare the active/visible tab, with the corrected range selected.
$ws.Activate() | Out-Null
$ws.Range("X2:Z2").Select() | Out-Null
